$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 286 and 287, pushing the existing rows
# 286-306 down to 288-308 (same formatting/styles carried over).
$ws.Range("A286:A287").EntireRow.Insert()

# Populate the first new row (286) - Zapallo / Camote / 1a (guarda)
$ws.Cells.Item(286, 1).Value = 3
$ws.Cells.Item(286, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(286, 3).Value = "Coquimbo"
$ws.Cells.Item(286, 4).Value = 44461
$ws.Cells.Item(286, 5).Value = 5
$ws.Cells.Item(286, 6).Value = 100112045
$ws.Cells.Item(286, 7).Value = "Zapallo"
$ws.Cells.Item(286, 8).Value = "Camote"
$ws.Cells.Item(286, 9).Value = "1a (guarda)"
$ws.Cells.Item(286, 10).Value = 90
$ws.Cells.Item(286, 11).Value = 800
$ws.Cells.Item(286, 12).Value = 800
$ws.Cells.Item(286, 13).Value = 800
$ws.Cells.Item(286, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(286, 15).Value = "Provincia de Talca"
$ws.Cells.Item(286, 16).Value = 800
$ws.Cells.Item(286, 17).Value = 1
$ws.Cells.Item(286, 18).Value = "Hortaliza"

# Populate the second new row (287) - Zapallo / Camote / 2a (guarda)
$ws.Cells.Item(287, 1).Value = 3
$ws.Cells.Item(287, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(287, 3).Value = "Coquimbo"
$ws.Cells.Item(287, 4).Value = 44461
$ws.Cells.Item(287, 5).Value = 5
$ws.Cells.Item(287, 6).Value = 100112045
$ws.Cells.Item(287, 7).Value = "Zapallo"
$ws.Cells.Item(287, 8).Value = "Camote"
$ws.Cells.Item(287, 9).Value = "2a (guarda)"
$ws.Cells.Item(287, 10).Value = 90
$ws.Cells.Item(287, 11).Value = 600
$ws.Cells.Item(287, 12).Value = 600
$ws.Cells.Item(287, 13).Value = 600
$ws.Cells.Item(287, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(287, 15).Value = "Provincia de Talca"
$ws.Cells.Item(287, 16).Value = 600
$ws.Cells.Item(287, 17).Value = 1
$ws.Cells.Item(287, 18).Value = "Hortaliza"
